$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 103. Excel shifts rows 103..172 down to 104..173,
# preserving each row's existing values/formatting.
$ws.Rows.Item(103).Insert()

# Populate the newly-inserted row 103 with a fresh weekly record. The fields
# that describe the series (market, region, product, variety, quality, unit,
# origin, classification, kg-or-unit factor) stay the same as the row that
# used to be here (now row 104); only the date and the volume/price figures
# change for this new week.
$ws.Range("A103").Value = 10
$ws.Range("B103").Value = "Vega Modelo de Temuco"
$ws.Range("C103").Value = "La Araucanía"
$ws.Range("D103").Value = 44582
$ws.Range("E103").Value = 9
$ws.Range("F103").Value = 100112005
$ws.Range("G103").Value = "Puerro"
$ws.Range("H103").Value = "Azul de Maquehue"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 40
$ws.Range("K103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 15000
$ws.Range("N103").Value = "$/docena de paquetes"
$ws.Range("O103").Value = "Provincia de Cautín"
$ws.Range("P103").Value = 1250
$ws.Range("Q103").Value = 12
$ws.Range("R103").Value = "Hortaliza"
